# Reorganizing project into three layers (command, Logic, revit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new project layout
$ws.Name = "revitSheet"

# Rename the Revit data table (name + displayName)
$lo = $ws.ListObjects.Item("RevitData")
$lo.Name = "revitData"

# Move the active selection on the sheet from J6 to B2
$ws.Range("B2").Select() | Out-Null
